$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B24").Value = 6371
$ws.Range("D24").Value = 5961189
$ws.Range("E24").Value = 935.6755611363993
$ws.Range("F24").Value = 8.608932833276505
$ws.Range("H24").Value = 26.28107356481375
